$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (sheet1) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value = 828
$ws.Cells.Item(7, 6).Value = 4233
$ws.Cells.Item(8, 3).Value = "北京·LookLook剧情式沉浸游戏互动动漫嘉年华（取消）"
$ws.Cells.Item(8, 4).Value = "东村文化创意产业园A1-2 五道杠实景片场"
$ws.Cells.Item(8, 5).Value = "2024.08.10 09:30-08.11 17:30"
$ws.Cells.Item(8, 6).Value = 1025
$ws.Cells.Item(8, 7).Value = "不可售"
$ws.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84741"
$ws.Cells.Item(8, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/WH4KKudj1716880619473.jpeg"
$ws.Cells.Item(9, 3).Value = "北京·不舍昼夜2.0-片羽拾光"
$ws.Cells.Item(9, 4).Value = "酒仙桥北路2号院798艺术区706后街1号 北京格瑞斯艺术酒店"
$ws.Cells.Item(9, 5).Value = "2024.08.10 10:30-08.11 02:00"
$ws.Cells.Item(9, 6).Value = 176
$ws.Cells.Item(9, 7).Value = 69
$ws.Cells.Item(9, 8).Value = "https://show.bilibili.com/platform/detail.html?id=88851"
$ws.Cells.Item(9, 9).Value = "//i1.hdslb.com/bfs/openplatform/202407/VACU64r21720427826079.png"
$ws.Cells.Item(10, 3).Value = "北京·广播剧《蝉女》专场活动"
$ws.Cells.Item(10, 4).Value = "北京展览馆 北京展览馆"
$ws.Cells.Item(10, 5).Value = "2024.08.10 11:50-08.10 15:10"
$ws.Cells.Item(10, 6).Value = 95
$ws.Cells.Item(10, 7).Value = 288
$ws.Cells.Item(10, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86947"
$ws.Cells.Item(10, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/ycrRjEPg1718176423186.jpeg"
$ws.Cells.Item(11, 3).Value = "北京·梦次元动漫展M30"
$ws.Cells.Item(11, 5).Value = "2024.08.10 10:00-08.11 17:00"
$ws.Cells.Item(11, 6).Value = 6157
$ws.Cells.Item(11, 7).Value = 80
$ws.Cells.Item(11, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83828"
$ws.Cells.Item(11, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/Qr2Bd5W41715931423636.jpeg"
$ws.Cells.Item(12, 6).Value = 6157
$ws.Cells.Item(15, 6).Value = 2354
$ws.Cells.Item(19, 6).Value = 9278
$ws.Cells.Item(21, 6).Value = 2505
$ws.Cells.Item(22, 6).Value = 196
$ws.Cells.Item(23, 6).Value = 2327
$ws.Cells.Item(24, 6).Value = 2472
$ws.Cells.Item(26, 6).Value = 246
$ws.Cells.Item(29, 6).Value = 61
$ws.Cells.Item(30, 6).Value = 335
$ws.Cells.Item(35, 6).Value = 74
$ws.Cells.Item(36, 6).Value = 385
$ws.Cells.Item(43, 6).Value = 2563

# ---- Sheet: 演出 (sheet2) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(22, 6).Value = 89

# ---- Sheet: 本地生活 (sheet3) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 696

# ---- Sheet: 全部类型 (sheet4) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 696
$ws.Cells.Item(5, 2).Value = "2024-07-17"
$ws.Cells.Item(5, 3).Value = "北京·“狐妖小红娘”限时快闪店"
$ws.Cells.Item(5, 4).Value = "王府井大街88号 北京王府井银泰in88购物中心"
$ws.Cells.Item(5, 5).Value = "2024.07.17 10:00-10.31 22:00"
$ws.Cells.Item(5, 6).Value = 104
$ws.Cells.Item(5, 7).Value = 98
$ws.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89613"
$ws.Cells.Item(5, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/n3TXriJX1721203778030.jpeg"
$ws.Cells.Item(6, 2).Value = "2024-07-20"
$ws.Cells.Item(6, 3).Value = "北京·英雄的苍穹：正子公也三国、水浒绘画艺术大展"
$ws.Cells.Item(6, 4).Value = "上庄大街18号附近郎园Park（地铁1号线八宝山b口） 郎园Park"
$ws.Cells.Item(6, 5).Value = "2024.07.20 10:00-08.18 19:00"
$ws.Cells.Item(6, 6).Value = 8
$ws.Cells.Item(6, 7).Value = 38
$ws.Cells.Item(6, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89424"
$ws.Cells.Item(6, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/hHCrntqE1721180587445.jpeg"
$ws.Cells.Item(7, 2).Value = "2024-07-21"
$ws.Cells.Item(7, 3).Value = "北京·航海王（ONE PIECE）25周年巡展"
$ws.Cells.Item(7, 4).Value = "酒仙桥路2号北京798艺术区A区 北京798艺术区"
$ws.Cells.Item(7, 5).Value = "2024.07.21 10:00-10.27 19:00"
$ws.Cells.Item(7, 6).Value = 369
$ws.Cells.Item(7, 7).Value = 98
$ws.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89233"
$ws.Cells.Item(7, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/WxL0mO9g1721011505489.png"
$ws.Cells.Item(8, 2).Value = "2024-08-02"
$ws.Cells.Item(8, 3).Value = "北京·沉浸悬疑剧 《穹顶弥声》"
$ws.Cells.Item(8, 4).Value = "通惠河畔1079金乾阁二层 光芒保利沉浸剧场"
$ws.Cells.Item(8, 5).Value = "2024.08.02 19:30-08.31 21:55"
$ws.Cells.Item(8, 6).Value = 5
$ws.Cells.Item(8, 7).Value = 128
$ws.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90182"
$ws.Cells.Item(8, 9).Value = "//i0.hdslb.com/bfs/openplatform/202406/KG9yeroD1719384539301.jpeg"
$ws.Cells.Item(9, 3).Value = "丰台·首家喜剧脱口秀魔仙喜剧 l 与您浪漫相约心动8月喜剧专场"
$ws.Cells.Item(9, 4).Value = "丽泽天地购物中心 丽泽天地购物中心"
$ws.Cells.Item(9, 5).Value = "2024.08.09 19:20-08.31 20:30"
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 59
$ws.Cells.Item(9, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90202"
$ws.Cells.Item(9, 9).Value = "//i1.hdslb.com/bfs/openplatform/202407/7A479TJD1722069826897.jpeg"
$ws.Cells.Item(11, 6).Value = 828
$ws.Cells.Item(12, 6).Value = 4233
$ws.Cells.Item(17, 3).Value = "北京·梦次元动漫展M30"
$ws.Cells.Item(17, 4).Value = "北京展览馆 北京展览馆"
$ws.Cells.Item(17, 5).Value = "2024.08.10 10:00-08.11 17:00"
$ws.Cells.Item(17, 6).Value = 6157
$ws.Cells.Item(17, 7).Value = 80
$ws.Cells.Item(17, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83828"
$ws.Cells.Item(17, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/Qr2Bd5W41715931423636.jpeg"
$ws.Cells.Item(18, 3).Value = "北京·狐妖小红娘专题聚会【免票活动】"
$ws.Cells.Item(18, 4).Value = "王府井大街88号 北京王府井银泰in88购物中心"
$ws.Cells.Item(18, 5).Value = "2024.08.10 14:00-08.10 18:00"
$ws.Cells.Item(18, 6).Value = 67
$ws.Cells.Item(18, 7).Value = 58
$ws.Cells.Item(18, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90238"
$ws.Cells.Item(18, 9).Value = "//i1.hdslb.com/bfs/openplatform/202408/mL8ytYCG1722578125040.jpeg"
$ws.Cells.Item(19, 3).Value = "北京·第五人格ONLY2.0"
$ws.Cells.Item(19, 4).Value = "永外高庄138号 北京大红门国际会展中心"
$ws.Cells.Item(19, 5).Value = "2024.08.10 10:00-08.10 17:00"
$ws.Cells.Item(19, 6).Value = 2354
$ws.Cells.Item(19, 7).Value = 70
$ws.Cells.Item(19, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86590"
$ws.Cells.Item(19, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/4jQBoo241716968548735.jpeg"
$ws.Cells.Item(22, 6).Value = 9278
$ws.Cells.Item(23, 2).Value = "2024-08-19"
$ws.Cells.Item(23, 3).Value = "北京·音阅派国漫演唱会-《一人之下》动画八周年专场演唱会"
$ws.Cells.Item(23, 4).Value = "中关村南大街33号中国国家图书馆内 国图艺术中心"
$ws.Cells.Item(23, 5).Value = "2024.08.19 19:30-08.19 21:00"
$ws.Cells.Item(23, 6).Value = 150
$ws.Cells.Item(23, 7).Value = 380
$ws.Cells.Item(23, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89052"
$ws.Cells.Item(23, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/wtYvGYL51720603864335.png"
$ws.Cells.Item(24, 3).Value = "北京·DICE CON 2024 第八届国际桌面游戏展"
$ws.Cells.Item(24, 4).Value = "东三环北路16 全国农业展览馆"
$ws.Cells.Item(24, 5).Value = "2024.08.23 13:00-08.25 18:00"
$ws.Cells.Item(24, 6).Value = 43
$ws.Cells.Item(24, 7).Value = 98
$ws.Cells.Item(24, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89579"
$ws.Cells.Item(24, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/sG2RY2Jk1721377035181.jpeg"
$ws.Cells.Item(25, 2).Value = "2024-08-23"
$ws.Cells.Item(25, 3).Value = "北京·喘气动漫嘉年华·暑期狂欢"
$ws.Cells.Item(25, 4).Value = "新风街3号 紫园·新风里"
$ws.Cells.Item(25, 5).Value = "2024.08.23 10:00-08.25 20:00"
$ws.Cells.Item(25, 6).Value = 2505
$ws.Cells.Item(25, 7).Value = 44.1
$ws.Cells.Item(25, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90022"
$ws.Cells.Item(25, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/Rb5sRto71722841653388.jpeg"
$ws.Cells.Item(26, 3).Value = " 北京·万游引力嘉年华 配音演员赵成晨&尘霜满眸 广播剧《奕曲同工》专场见面&签售会"
$ws.Cells.Item(26, 4).Value = "金蝉西路甲1号（地铁七号线南楼梓庄站） 北京酷车国际汇展中心"
$ws.Cells.Item(26, 5).Value = "2024.08.24 11:00-08.24 17:00"
$ws.Cells.Item(26, 6).Value = 196
$ws.Cells.Item(26, 7).Value = 288
$ws.Cells.Item(26, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89054"
$ws.Cells.Item(26, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/FadWpN3x1720599868028.jpeg"
$ws.Cells.Item(27, 6).Value = 2472
$ws.Cells.Item(29, 6).Value = 246
$ws.Cells.Item(32, 6).Value = 61
$ws.Cells.Item(33, 6).Value = 335
$ws.Cells.Item(34, 2).Value = "2024-09-07"
$ws.Cells.Item(34, 3).Value = "北京·明日方舟ONLY同人展"
$ws.Cells.Item(34, 4).Value = "魏永路9-1号 中国书画院"
$ws.Cells.Item(34, 5).Value = "2024.09.07 10:00-09.08 18:00"
$ws.Cells.Item(34, 6).Value = 45
$ws.Cells.Item(34, 7).Value = 65
$ws.Cells.Item(34, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90218"
$ws.Cells.Item(34, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/h3axTTjv1722998734230.jpeg"
$ws.Cells.Item(36, 3).Value = "北京·AP动漫嘉年华"
$ws.Cells.Item(36, 4).Value = "永外高庄138号 北京大红门国际会展中心"
$ws.Cells.Item(36, 5).Value = "2024.09.15 10:00-09.15 17:00"
$ws.Cells.Item(36, 6).Value = 44
$ws.Cells.Item(36, 7).Value = 60
$ws.Cells.Item(36, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89165"
$ws.Cells.Item(36, 9).Value = "//i1.hdslb.com/bfs/openplatform/202407/vyuPGUrJ1720747471465.jpeg"
$ws.Cells.Item(37, 3).Value = "北京·ICOS SP漫展04动漫节"
$ws.Cells.Item(37, 4).Value = "金蝉西路甲1号（地铁七号线南楼梓庄站） 北京酷车国际汇展中心"
$ws.Cells.Item(37, 5).Value = "2024.09.15 09:00-09.16 17:00"
$ws.Cells.Item(37, 6).Value = 74
$ws.Cells.Item(37, 7).Value = 80
$ws.Cells.Item(37, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90286"
$ws.Cells.Item(37, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/tPazRaBV1722595834650.jpeg"
$ws.Cells.Item(38, 2).Value = "2024-09-15"
$ws.Cells.Item(38, 3).Value = "北京·MQ&THEBONE首届怀旧同人only"
$ws.Cells.Item(38, 4).Value = "安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园L1层"
$ws.Cells.Item(38, 5).Value = "2024.09.15 10:00-09.16 17:00"
$ws.Cells.Item(38, 6).Value = 385
$ws.Cells.Item(38, 7).Value = 6.6
$ws.Cells.Item(38, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90096"
$ws.Cells.Item(38, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/BBmePlWJ1722409048209.jpeg"
$ws.Cells.Item(39, 2).Value = "2024-09-15"
$ws.Cells.Item(39, 3).Value = "北京·原神only4.0同人展"
$ws.Cells.Item(39, 4).Value = "北花园路1号 超级蜂巢"
$ws.Cells.Item(39, 5).Value = "2024.09.15 10:00-09.15 17:00"
$ws.Cells.Item(39, 6).Value = 1224
$ws.Cells.Item(39, 7).Value = 68
$ws.Cells.Item(39, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87564"
$ws.Cells.Item(39, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/EfEAeJDS1720776874376.jpeg"
$ws.Cells.Item(40, 2).Value = "2024-09-16"
$ws.Cells.Item(40, 3).Value = "北京·原神×星穹铁道only2.0同人展"
$ws.Cells.Item(40, 4).Value = "高碑店东路超级蜂巢 5G直播基地"
$ws.Cells.Item(40, 5).Value = "2024.09.16 10:00-09.16 17:00"
$ws.Cells.Item(40, 6).Value = 1223
$ws.Cells.Item(40, 7).Value = 68
$ws.Cells.Item(40, 8).Value = "https://show.bilibili.com/platform/detail.html?id=88285"
$ws.Cells.Item(40, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/iWlE3Q9X1719554169582.jpeg"
$ws.Cells.Item(41, 2).Value = "2024-09-17"
$ws.Cells.Item(41, 3).Value = "北京·双男主only之皎皎秋月夜"
$ws.Cells.Item(41, 4).Value = "太平庄中街西端 北京天通苑黄河京都会议中心"
$ws.Cells.Item(41, 5).Value = "2024.09.17 10:00-09.17 17:00"
$ws.Cells.Item(41, 6).Value = 75
$ws.Cells.Item(41, 7).Value = 79
$ws.Cells.Item(41, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89763"
$ws.Cells.Item(41, 9).Value = "//i1.hdslb.com/bfs/openplatform/202407/nUiFpHBb1721723099117.jpeg"
$ws.Cells.Item(42, 2).Value = "2024-09-17"
$ws.Cells.Item(42, 3).Value = "北京·马娘ONLY2"
$ws.Cells.Item(42, 4).Value = "永外高庄138号 北京大红门国际会展中心"
$ws.Cells.Item(42, 5).Value = "2024.09.17 10:00-09.17 17:00"
$ws.Cells.Item(42, 6).Value = 101
$ws.Cells.Item(42, 7).Value = 75
$ws.Cells.Item(42, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89334"
$ws.Cells.Item(42, 9).Value = "//i1.hdslb.com/bfs/openplatform/202407/0LLxCfKo1721112673092.png"
$ws.Cells.Item(43, 6).Value = 2563
$ws.Cells.Item(46, 2).Value = "2024-10-25"
$ws.Cells.Item(46, 3).Value = "北京·伦敦西区音乐剧明星演唱会-经典版"
$ws.Cells.Item(46, 5).Value = "2024.10.25 19:30-10.26 21:30"
$ws.Cells.Item(46, 6).Value = 3
$ws.Cells.Item(46, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89359"
$ws.Cells.Item(46, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/PzPiEKUI1721114840552.jpeg"
$ws.Cells.Item(47, 3).Value = "北京·伦敦西区音乐剧明星演唱会（摇滚版）"
$ws.Cells.Item(47, 4).Value = "西直门外大街135号（北京展览馆内） 北京展览馆剧场"
$ws.Cells.Item(47, 5).Value = "2024.10.26 14:30-10.26 16:30"
$ws.Cells.Item(47, 6).Value = 4
$ws.Cells.Item(47, 7).Value = 144
$ws.Cells.Item(47, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89400"
$ws.Cells.Item(47, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/TYPRpfu21721116217467.jpeg"
$ws.Cells.Item(48, 2).Value = "2024-11-09"
$ws.Cells.Item(48, 3).Value = "北京·Aw动漫游戏嘉年华9th"
$ws.Cells.Item(48, 4).Value = "石景山路68号 北京首钢会展中心"
$ws.Cells.Item(48, 5).Value = "2024.11.09 09:30-11.10 17:30"
$ws.Cells.Item(48, 6).Value = 27
$ws.Cells.Item(48, 7).Value = 60
$ws.Cells.Item(48, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90357"
$ws.Cells.Item(48, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/Nl370YWS1722588270723.jpeg"
$ws.Cells.Item(50, 6).Value = 89
$ws.Cells.Item(51, 6).Value = 89

Write-Host "Applied all cell updates."